# Formatação do Escopo, Requisitos e User Stories.
#
# This script:
#  1. Splits the run "Problematização" into "Problematiza" + "ção" and
#     relocates the "_GoBack" bookmark into the gap (moving it away from
#     the "Não escopo" heading later in the document).
#  2. Justifies (wdAlignParagraphJustify = 3) every paragraph from
#     "PORTAL ECONOMUNDI" through the end of the document body, including
#     the two empty paragraphs that used to be bare <w:p/> elements, and
#     every paragraph inside the stakeholders-approval table.

$d = $word.ActiveDocument
$wdAlignParagraphJustify = 3

# ---------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark to inside "Problematização", splitting
#    the run into "Problematiza" | (bookmark) | "ção".
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$paraStart = $p5.Range.Start
$paraText = $p5.Range.Text
$splitOffset = $paraText.IndexOf("ção", $paraText.IndexOf(" e Problematiza"))
$splitPos = $paraStart + $splitOffset
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------
# 2. Justify every paragraph in the main body from "PORTAL ECONOMUNDI"
#    (paragraph 4) through the last regular body paragraph before the
#    approval table (paragraph 90, ending in "Data").
# ---------------------------------------------------------------------
for ($i = 4; $i -le 90; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
}

# ---------------------------------------------------------------------
# 3. Justify every paragraph inside the approval table (both rows).
# ---------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
    }
}

# ---------------------------------------------------------------------
# 4. Justify the two trailing empty paragraphs after the table.
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
$tailFirst = $d.Paragraphs.Item($total - 1)
$tailFirst.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
$tailLast = $d.Paragraphs.Item($total)
$tailLast.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
